$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("weight"), shifting old D (grade) -> E and old E (infected) -> F
$ws.Columns("D").Insert()

# Header for the new column
$ws.Range("D1").Value = "weight"

# Corrected jump_length values for rows 5-7 (column B)
$ws.Range("B5").Value = 5.6
$ws.Range("B6").Value = 9.1
$ws.Range("B7").Value = 8.2

# New weight values (column D) for each data row
$weights = @{
    2  = 2.1
    3  = 2.3
    4  = 2.8
    5  = 2.4
    6  = 1.2
    7  = 4.1
    8  = 3.2
    9  = 1.1
    10 = 2.1
    11 = 2.4
    12 = 2.1
    13 = 1.5
    14 = 3.7
    15 = 2.9
    16 = 3.1
    17 = 4.2
    18 = 5.1
    19 = 3.5
    20 = 3.2
    21 = 4.6
    22 = 3.7
}

foreach ($row in $weights.Keys) {
    $ws.Range("D$row").Value = $weights[$row]
}
